# Updates Leve profit-calculation cells (currentAveragePrice* / LevePrice* /
# LeveProfit* columns, H:N) across several job sheets, per scheduled-runner
# price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 76
$ws.Range("H76").Value = 3973.45
$ws.Range("I76").Value = 3250.25
$ws.Range("J76").Value = 4455.5835
$ws.Range("K76").Value = 3250.25
$ws.Range("L76").Value = 4455.5835
$ws.Range("M76").Value = -2935.25
$ws.Range("N76").Value = -5085.5835

# row 79
$ws.Range("H79").Value = 3973.45
$ws.Range("I79").Value = 3250.25
$ws.Range("J79").Value = 4455.5835
$ws.Range("K79").Value = 3250.25
$ws.Range("L79").Value = 4455.5835
$ws.Range("M79").Value = -2158.25
$ws.Range("N79").Value = -6639.5835

# row 108
$ws.Range("H108").Value = 24833
$ws.Range("J108").Value = 24833
$ws.Range("L108").Value = 24833
$ws.Range("N108").Value = -32513

# row 110
$ws.Range("H110").Value = 32997.5
$ws.Range("J110").Value = 32997.5
$ws.Range("L110").Value = 32997.5
$ws.Range("N110").Value = -41177.5

# row 111
$ws.Range("H111").Value = 9520.733
$ws.Range("I111").Value = 15753.857
$ws.Range("J111").Value = 4066.75
$ws.Range("K111").Value = 47261.571
$ws.Range("L111").Value = 12200.25
$ws.Range("M111").Value = -44194.571
$ws.Range("N111").Value = -18334.25

# row 123
$ws.Range("H123").Value = 33101.6
$ws.Range("J123").Value = 33101.6
$ws.Range("L123").Value = 33101.6
$ws.Range("N123").Value = -42901.6

$ws = $wb.Worksheets.Item("ARM")
# row 26
$ws.Range("H26").Value = 4449.625
$ws.Range("I26").Value = 3656.7144
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 3656.7144
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = -3326.7144
$ws.Range("N26").Value = -10660

# row 98
$ws.Range("H98").Value = 15407
$ws.Range("J98").Value = 15407
$ws.Range("L98").Value = 15407
$ws.Range("N98").Value = -21397

$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 74129.71000000001
$ws.Range("I20").Value = 79562.766
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 79562.766
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -79315.766
$ws.Range("N20").Value = -3994

# row 107
$ws.Range("H107").Value = 55603828
$ws.Range("I107").Value = 66724336
$ws.Range("J107").Value = 1295.6666
$ws.Range("K107").Value = 66724336
$ws.Range("L107").Value = 1295.6666
$ws.Range("M107").Value = -66722416
$ws.Range("N107").Value = -5135.6666

$ws = $wb.Worksheets.Item("CRP")
# row 105
$ws.Range("H105").Value = 1807.1904
$ws.Range("I105").Value = 1742.8823
$ws.Range("J105").Value = 2080.5
$ws.Range("K105").Value = 1742.8823
$ws.Range("L105").Value = 2080.5
$ws.Range("M105").Value = 4.117700000000013
$ws.Range("N105").Value = -5574.5

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1329.7941
$ws.Range("I5").Value = 596.913
$ws.Range("J5").Value = 2862.182
$ws.Range("K5").Value = 1790.739
$ws.Range("L5").Value = 8586.545999999998
$ws.Range("M5").Value = -1678.739
$ws.Range("N5").Value = -8810.545999999998

# row 114
$ws.Range("H114").Value = 604.9286
$ws.Range("I114").Value = 241
$ws.Range("J114").Value = 1090.1666
$ws.Range("K114").Value = 723
$ws.Range("L114").Value = 3270.4998
$ws.Range("M114").Value = 2531
$ws.Range("N114").Value = -9778.4998

# row 135
$ws.Range("H135").Value = 1329.7941
$ws.Range("I135").Value = 596.913
$ws.Range("J135").Value = 2862.182
$ws.Range("K135").Value = 5372.217000000001
$ws.Range("L135").Value = 25759.638
$ws.Range("M135").Value = -2837.217000000001
$ws.Range("N135").Value = -30829.638

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 96150.45
$ws.Range("I70").Value = 203711.2
$ws.Range("J70").Value = 6516.5
$ws.Range("K70").Value = 203711.2
$ws.Range("L70").Value = 6516.5
$ws.Range("M70").Value = -203441.2
$ws.Range("N70").Value = -7056.5

# row 73
$ws.Range("H73").Value = 96150.45
$ws.Range("I73").Value = 203711.2
$ws.Range("J73").Value = 6516.5
$ws.Range("K73").Value = 203711.2
$ws.Range("L73").Value = 6516.5
$ws.Range("M73").Value = -202775.2
$ws.Range("N73").Value = -8388.5

# row 80
$ws.Range("H80").Value = 250002220
$ws.Range("I80").Value = 500002500
$ws.Range("J80").Value = 1973
$ws.Range("K80").Value = 500002500
$ws.Range("L80").Value = 1973
$ws.Range("M80").Value = -500001502
$ws.Range("N80").Value = -3969

# row 83
$ws.Range("H83").Value = 250002220
$ws.Range("I83").Value = 500002500
$ws.Range("J83").Value = 1973
$ws.Range("K83").Value = 2500012500
$ws.Range("L83").Value = 9865
$ws.Range("M83").Value = -2500007508
$ws.Range("N83").Value = -19849

# row 113
$ws.Range("H113").Value = 2569.2144
$ws.Range("I113").Value = 3664.6
$ws.Range("J113").Value = 1960.6666
$ws.Range("K113").Value = 3664.6
$ws.Range("L113").Value = 1960.6666
$ws.Range("M113").Value = -1494.6
$ws.Range("N113").Value = -6300.6666

# row 132
$ws.Range("H132").Value = 2758.4092
$ws.Range("I132").Value = 1736.6875
$ws.Range("J132").Value = 5483
$ws.Range("K132").Value = 5210.0625
$ws.Range("L132").Value = 16449
$ws.Range("M132").Value = -2680.0625
$ws.Range("N132").Value = -21509

$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 2679.55
$ws.Range("I61").Value = 2079
$ws.Range("J61").Value = 3170.9092
$ws.Range("K61").Value = 2079
$ws.Range("L61").Value = 3170.9092
$ws.Range("M61").Value = -1877
$ws.Range("N61").Value = -3574.9092

# row 68
$ws.Range("H68").Value = 2959.5
$ws.Range("I68").Value = 1859.9
$ws.Range("K68").Value = 1859.9
$ws.Range("M68").Value = -1110.9

# row 71
$ws.Range("H71").Value = 2959.5
$ws.Range("I71").Value = 1859.9
$ws.Range("K71").Value = 9299.5
$ws.Range("M71").Value = -5555.5

# row 113
$ws.Range("H113").Value = 2679.55
$ws.Range("I113").Value = 2079
$ws.Range("J113").Value = 3170.9092
$ws.Range("K113").Value = 2079
$ws.Range("L113").Value = 3170.9092
$ws.Range("M113").Value = 91
$ws.Range("N113").Value = -7510.9092

$ws = $wb.Worksheets.Item("WVR")
# row 45
$ws.Range("H45").Value = 7258.1665
$ws.Range("J45").Value = 7258.1665
$ws.Range("L45").Value = 7258.1665
$ws.Range("N45").Value = -8240.166499999999

# row 81
$ws.Range("H81").Value = 134204.8
$ws.Range("I81").Value = 101006.8
$ws.Range("J81").Value = 200600.8
$ws.Range("K81").Value = 202013.6
$ws.Range("L81").Value = 401201.6
$ws.Range("M81").Value = -200952.6
$ws.Range("N81").Value = -403323.6

# row 84
$ws.Range("H84").Value = 134204.8
$ws.Range("I84").Value = 101006.8
$ws.Range("J84").Value = 200600.8
$ws.Range("K84").Value = 1010068
$ws.Range("L84").Value = 2006008
$ws.Range("M84").Value = -1004764
$ws.Range("N84").Value = -2016616

# row 117
$ws.Range("H117").Value = 32828.57
$ws.Range("J117").Value = 32828.57
$ws.Range("L117").Value = 32828.57
$ws.Range("N117").Value = -42006.57

# row 124
$ws.Range("H124").Value = 45995
$ws.Range("J124").Value = 45995
$ws.Range("L124").Value = 45995
$ws.Range("N124").Value = -55815
